$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): relabel the parameter-list / group-desc columns ---
# PARAM_LIST -> PARAM_ND_LIST (col F header)
# (new) PARAM_HD_LIST inserted as col G header
# PARAM_GRP_DESC moves to col H header
$ws.Range("F1").Value = "PARAM_ND_LIST"
$ws.Range("G1").Value = "PARAM_HD_LIST"
$ws.Range("H1").Value = "PARAM_GRP_DESC"

# --- Data rows: column F (previously V1000_LIST_LV/HV) -> V1000_*_ND ---
# --- column G (previously V1000_*_HD_MODS) -> V1000_*_HD ---
# --- column H (previously V1000_GRP_DESC) stays V1000_GRP_DESC ---
$lv_nd = "V1000_LV_ND"
$hv_nd = "V1000_HV_ND"
$lv_hd = "V1000_LV_HD"
$hv_hd = "V1000_HV_HD"

# Rows 2,4,6,8 are the "LV" (low voltage) drive rows; rows 3,5,7,9 are "HV" rows
$ws.Range("F2").Value = $lv_nd
$ws.Range("G2").Value = $lv_hd
$ws.Range("H2").Value = "V1000_GRP_DESC"

$ws.Range("F3").Value = $hv_nd
$ws.Range("G3").Value = $hv_hd
$ws.Range("H3").Value = "V1000_GRP_DESC"

$ws.Range("F4").Value = $lv_nd
$ws.Range("G4").Value = $lv_hd
$ws.Range("H4").Value = "V1000_GRP_DESC"

$ws.Range("F5").Value = $hv_nd
$ws.Range("G5").Value = $hv_hd
$ws.Range("H5").Value = "V1000_GRP_DESC"

$ws.Range("F6").Value = $lv_nd
$ws.Range("G6").Value = $lv_hd
$ws.Range("H6").Value = "V1000_GRP_DESC"

$ws.Range("F7").Value = $hv_nd
$ws.Range("G7").Value = $hv_hd
$ws.Range("H7").Value = "V1000_GRP_DESC"

$ws.Range("F8").Value = $lv_nd
$ws.Range("G8").Value = $lv_hd
$ws.Range("H8").Value = "V1000_GRP_DESC"

$ws.Range("F9").Value = $hv_nd
$ws.Range("G9").Value = $hv_hd
$ws.Range("H9").Value = "V1000_GRP_DESC"

# --- Widen column F to fit the new, longer header text ---
# (target stored width is 15.5546875 chars; engine rounds ColumnWidth to whole
#  pixels using a 7px/char metric, so 14.86 is the closest input that lands on
#  the nearest achievable pixel boundary)
$ws.Columns.Item(6).ColumnWidth = 14.86

# --- Move the active selection to I8 ---
$ws.Range("I8").Select()

# --- Shift the saved workbook window position (xWindow); best effort, the
#     window-geometry attributes are cosmetic/not part of the COM object model ---
$excel.Left = 6615
$excel.ActiveWindow.Left = 6615
